# Update a small number of trial values in the L1/L2 sheets to reflect
# the longer epoch lengths (slight LI value changes: 1 -> -1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("L1")
$ws2 = $wb.Worksheets.Item("L2")

# Sheet L1 changes
$ws1.Range("V5").Value = -1
$ws1.Range("O12").Value = -1

# Sheet L2 changes
$ws2.Range("C2").Value = -1
$ws2.Range("H4").Value = -1
$ws2.Range("H5").Value = -1
$ws2.Range("K7").Value = -1
$ws2.Range("O7").Value = -1
$ws2.Range("F8").Value = -1
$ws2.Range("M9").Value = -1
$ws2.Range("T10").Value = -1
$ws2.Range("I16").Value = -1
$ws2.Range("S20").Value = -1
$ws2.Range("U21").Value = -1
$ws2.Range("B23").Value = -1
$ws2.Range("B24").Value = -1
$ws2.Range("B26").Value = -1
